# Updated files from WRI China bringing the Hong Kong EPS up to v2.0.0
#
# Changes applied to "Pot Perc Red in Fuel Use fr Inc Cogen and WHR.xlsx":
#   1. Remove the extra, empty "Sheet1" worksheet that sat between "About"
#      and "Data".
#   2. On the "About" sheet, clear out the stray "use US value" note that
#      lived in B13 (below the Notes section).
#   3. On the "PPRiFUfICaWHR" sheet, rename the header in B1 from
#      "Pot Perc Red in Fuel Use" to "Pot Perc Red in Fuel Use
#      (dimensionless)", wrapping the text and growing row 1 to fit it.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# 1. Drop the unused "Sheet1" tab entirely.
$extraSheet = $wb.Worksheets.Item("Sheet1")
$extraSheet.Delete()

# 2. "About" sheet: remove the leftover red "use US value" annotation.
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("B13").Clear()

# 3. "PPRiFUfICaWHR" sheet: update the column header text/formatting.
$wsPPR = $wb.Worksheets.Item("PPRiFUfICaWHR")
$wsPPR.Range("B1").Value = "Pot Perc Red in Fuel Use (dimensionless)"
$wsPPR.Range("B1").WrapText = $true
$wsPPR.Rows.Item(1).RowHeight = 28.5

# Restore sensible selections on the affected sheets.
$wsPPR.Activate()
$wsPPR.Range("B1").Select()

$wsAbout.Activate()
$wsAbout.Range("A1").Select()
